$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 33334672
$ws.Range("J32").Value = 1434.875
$ws.Range("L32").Value = 1434.875
$ws.Range("N32").Value = -2086.875
$ws.Range("H51").Value = 5298168.5
$ws.Range("J51").Value = 7521.1
$ws.Range("L51").Value = 7521.1
$ws.Range("N51").Value = -8489.1
$ws.Range("H76").Value = 3069.7585
$ws.Range("I76").Value = 3000.1304
$ws.Range("J76").Value = 3336.6667
$ws.Range("K76").Value = 3000.1304
$ws.Range("L76").Value = 3336.6667
$ws.Range("M76").Value = -2685.1304
$ws.Range("N76").Value = -3966.6667
$ws.Range("H79").Value = 3069.7585
$ws.Range("I79").Value = 3000.1304
$ws.Range("J79").Value = 3336.6667
$ws.Range("K79").Value = 3000.1304
$ws.Range("L79").Value = 3336.6667
$ws.Range("M79").Value = -1908.1304
$ws.Range("N79").Value = -5520.6667
$ws.Range("H105").Value = 48887.668
$ws.Range("J105").Value = 48887.668
$ws.Range("L105").Value = 48887.668
$ws.Range("N105").Value = -55875.668
$ws.Range("H112").Value = 1901.8
$ws.Range("J112").Value = 2046
$ws.Range("L112").Value = 6138
$ws.Range("N112").Value = -8354
$ws.Range("H129").Value = 1361.7878
$ws.Range("I129").Value = 2532.6
$ws.Range("J129").Value = 1152.7142
$ws.Range("K129").Value = 7597.799999999999
$ws.Range("L129").Value = 3458.1426
$ws.Range("M129").Value = -2597.799999999999
$ws.Range("N129").Value = -13458.1426
$ws.Range("H132").Value = 2948.93
$ws.Range("I132").Value = 1640.241
$ws.Range("J132").Value = 9338.412
$ws.Range("K132").Value = 4920.723
$ws.Range("L132").Value = 28015.236
$ws.Range("M132").Value = -2390.723
$ws.Range("N132").Value = -33075.236
$ws.Range("H137").Value = 4482.5
$ws.Range("I137").Value = 1500
$ws.Range("J137").Value = 6470.8335
$ws.Range("K137").Value = 4500
$ws.Range("L137").Value = 19412.5005
$ws.Range("M137").Value = -1950
$ws.Range("N137").Value = -24512.5005
$ws.Range("H141").Value = 1979.3125
$ws.Range("I141").Value = 799.9268
$ws.Range("J141").Value = 8887.143
$ws.Range("K141").Value = 2399.7804
$ws.Range("L141").Value = 26661.429
$ws.Range("M141").Value = 2780.2196
$ws.Range("N141").Value = -37021.429
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9234.671
$ws.Range("I32").Value = 8830.950999999999
$ws.Range("J32").Value = 11510.182
$ws.Range("K32").Value = 8830.950999999999
$ws.Range("L32").Value = 11510.182
$ws.Range("M32").Value = -8543.950999999999
$ws.Range("N32").Value = -12084.182
$ws.Range("H80").Value = 45628.223
$ws.Range("J80").Value = 45628.223
$ws.Range("L80").Value = 45628.223
$ws.Range("N80").Value = -47624.223
$ws.Range("H83").Value = 45628.223
$ws.Range("J83").Value = 45628.223
$ws.Range("L83").Value = 136884.669
$ws.Range("N83").Value = -146868.669
$ws.Range("H88").Value = 15488831
$ws.Range("I88").Value = 66669100
$ws.Range("J88").Value = 2693764.2
$ws.Range("K88").Value = 66669100
$ws.Range("L88").Value = 2693764.2
$ws.Range("M88").Value = -66668694
$ws.Range("N88").Value = -2694576.2
$ws.Range("H91").Value = 15488831
$ws.Range("I91").Value = 66669100
$ws.Range("J91").Value = 2693764.2
$ws.Range("K91").Value = 66669100
$ws.Range("L91").Value = 2693764.2
$ws.Range("M91").Value = -66667696
$ws.Range("N91").Value = -2696572.2
$ws.Range("H102").Value = 11725.182
$ws.Range("I102").Value = 1940.2727
$ws.Range("J102").Value = 21510.092
$ws.Range("K102").Value = 1940.2727
$ws.Range("L102").Value = 21510.092
$ws.Range("M102").Value = -318.2727
$ws.Range("N102").Value = -24754.092
$ws.Range("H138").Value = 43500
$ws.Range("J138").Value = 43500
$ws.Range("L138").Value = 43500
$ws.Range("N138").Value = -53780
$ws.Range("H139").Value = 44387.668
$ws.Range("J139").Value = 44387.668
$ws.Range("L139").Value = 44387.668
$ws.Range("N139").Value = -54667.668
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 5808.6665
$ws.Range("I54").Value = 5808.6665
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 5808.6665
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -5324.6665
$ws.Range("N54").Value = ""
$ws.Range("H134").Value = 3012.7974
$ws.Range("I134").Value = 2436.4583
$ws.Range("K134").Value = 7309.374899999999
$ws.Range("M134").Value = -4774.374899999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 774.2
$ws.Range("I16").Value = 734.1667
$ws.Range("K16").Value = 734.1667
$ws.Range("M16").Value = -447.1667
$ws.Range("H31").Value = 2862.42
$ws.Range("I31").Value = 1156.6552
$ws.Range("J31").Value = 3559.1409
$ws.Range("K31").Value = 1156.6552
$ws.Range("L31").Value = 3559.1409
$ws.Range("M31").Value = -861.6551999999999
$ws.Range("N31").Value = -4149.1409
$ws.Range("H34").Value = 2862.42
$ws.Range("I34").Value = 1156.6552
$ws.Range("J34").Value = 3559.1409
$ws.Range("K34").Value = 1156.6552
$ws.Range("L34").Value = 3559.1409
$ws.Range("M34").Value = -954.6551999999999
$ws.Range("N34").Value = -3963.1409
$ws.Range("H54").Value = 12546
$ws.Range("J54").Value = 12546
$ws.Range("L54").Value = 12546
$ws.Range("N54").Value = -13862
$ws.Range("H113").Value = 774.2
$ws.Range("I113").Value = 734.1667
$ws.Range("K113").Value = 734.1667
$ws.Range("M113").Value = 1435.8333
$ws.Range("H132").Value = 43599.44
$ws.Range("I132").Value = 1446.32
$ws.Range("K132").Value = 4338.96
$ws.Range("M132").Value = -1808.96
$ws.Range("H134").Value = 638322.5
$ws.Range("I134").Value = 1376.0769
$ws.Range("K134").Value = 4128.2307
$ws.Range("M134").Value = -1593.2307
$ws.Range("H137").Value = 37599.668
$ws.Range("J137").Value = 37599.668
$ws.Range("L137").Value = 37599.668
$ws.Range("N137").Value = -47799.668
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 138969.06
$ws.Range("I56").Value = 138969.06
$ws.Range("K56").Value = 138969.06
$ws.Range("M56").Value = -138439.06
$ws.Range("H64").Value = 2579.8
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 2579.8
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 7739.400000000001
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = -8279.400000000001
$ws.Range("H67").Value = 2579.8
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2579.8
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 7739.400000000001
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = -9611.400000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 23795
$ws.Range("J6").Value = 23795
$ws.Range("L6").Value = 23795
$ws.Range("N6").Value = -24021
$ws.Range("H16").Value = 23795
$ws.Range("J16").Value = 23795
$ws.Range("L16").Value = 23795
$ws.Range("N16").Value = -24295
$ws.Range("H138").Value = 54000
$ws.Range("J138").Value = 54000
$ws.Range("L138").Value = 54000
$ws.Range("N138").Value = -64280
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 70001.5
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("H26").Value = 6000
$ws.Range("J26").Value = 6000
$ws.Range("L26").Value = 6000
$ws.Range("N26").Value = -6590
$ws.Range("H40").Value = 5225.5713
$ws.Range("I40").Value = 3118.8
$ws.Range("K40").Value = 3118.8
$ws.Range("M40").Value = -2982.8
$ws.Range("H122").Value = 2444.4443
$ws.Range("I122").Value = 2357.1428
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 7071.428400000001
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -4621.428400000001
$ws.Range("N122").Value = -13150
$ws.Range("H134").Value = 50999.668
$ws.Range("J134").Value = 50999.668
$ws.Range("L134").Value = 50999.668
$ws.Range("N134").Value = -61139.668
$ws.Range("H136").Value = 1213.1569
$ws.Range("I136").Value = 974.95123
$ws.Range("K136").Value = 2924.85369
$ws.Range("M136").Value = -374.8536899999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 33938
$ws.Range("J26").Value = 33938
$ws.Range("L26").Value = 33938
$ws.Range("N26").Value = -34524
$ws.Range("H100").Value = 436
$ws.Range("I100").Value = 436
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 872
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -331
$ws.Range("N100").Value = ""
$ws.Range("H107").Value = 6667593.5
$ws.Range("I107").Value = 989.1111
$ws.Range("J107").Value = 16667500
$ws.Range("K107").Value = 2967.3333
$ws.Range("L107").Value = 50002500
$ws.Range("M107").Value = -1047.3333
$ws.Range("N107").Value = -50006340
$ws.Range("H122").Value = 2598063.5
$ws.Range("I122").Value = 2857769.8
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 8573309.399999999
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -8570859.399999999
$ws.Range("N122").Value = -7900
$ws.Range("H131").Value = 50121
$ws.Range("J131").Value = 50121
$ws.Range("L131").Value = 50121
$ws.Range("N131").Value = -60201
$ws.Range("H132").Value = 1659.6377
$ws.Range("I132").Value = 1378.8302
$ws.Range("K132").Value = 4136.4906
$ws.Range("M132").Value = -1606.4906
$ws.Range("H139").Value = 57550
$ws.Range("J139").Value = 57550
$ws.Range("L139").Value = 57550
$ws.Range("N139").Value = -67830
